# test_create_one_column.xlsx - create_one_onglet_by_participant was changed
# to write its per-module Oui/Non answer columns onto this workbook too, so
# the "sheet1" export now carries two extra repeats of the four-column
# Alain/Henri/Tony/Dulcinée block (8 columns) right before the trailing
# "Adresse de courriel" / empty column pair.
#
# Concretely: insert 8 new columns at AFM:AFT (pushing the old AFM/AFN pair
# - the e-mail address and the empty numeric cell - to AFU/AFV) and fill the
# new columns with a copy of the existing I:P block (same values, same
# style) for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Make room for the new 8-column block; existing columns from AFM onward
# (email address + empty cell) shift right to AFU/AFV automatically.
$ws.Columns("AFM:AFT").Insert()

# Populate the freshly inserted columns with the same repeating
# Alain/Henri/Tony/Dulcinée (header row) / OUI-NON (data rows) pattern
# already used by columns I:P, preserving their style.
$ws.Range("I1:P9").Copy()
$ws.Range("AFM1:AFT9").PasteSpecial()

$excel.CutCopyMode = $false
